# Workbook / worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Extend the header row formatting (bold / centered / bordered style
#    that already lives on A1:C1) across the new columns D1:G1.
# ---------------------------------------------------------------------
$ws.Range("A1:C1").Copy()
$ws.Range("D1:G1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. The date/time formatting used to live on column A (rows 2-5).
#    In the new layout the date lives in column B, so copy that
#    number format over to B2:B13 before we overwrite column A.
# ---------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("B2:B13").PasteSpecial(-4122)

# Column A no longer holds dates - restore it to the default style.
$ws.Range("A2:A13").Style = "Normal"

# ---------------------------------------------------------------------
# 3. Header row values / new column headers.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Transaction_ID"
$ws.Range("B1").Value = "Date"
$ws.Range("C1").Value = "Total_Price"
$ws.Range("D1").Value = "Unit_Price"
$ws.Range("E1").Value = "Quantity"
$ws.Range("F1").Value = "Product_ID"
$ws.Range("G1").Value = "Product_Name"

# ---------------------------------------------------------------------
# 4. Data rows.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = 45233.65696068287
$ws.Range("C2").Value = 4
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = "Margerita"

$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 45233.65696068287
$ws.Range("C3").Value = 8
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = "Spicy"

$ws.Range("A4").Value = 5
$ws.Range("B4").Value = 45233.65696068287
$ws.Range("C4").Value = 7
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = "Fancy"

$ws.Range("A5").Value = 6
$ws.Range("B5").Value = 45233.65892061342
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = "Margerita"

$ws.Range("A6").Value = 6
$ws.Range("B6").Value = 45233.65892061342
$ws.Range("C6").Value = 12
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = "Spicy"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 45233.65892061342
$ws.Range("C7").Value = 28
$ws.Range("D7").Value = 7
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = "Fancy"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 45233.65996528935
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = "Margerita"

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 45233.65996528935
$ws.Range("C9").Value = 20
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = "Spicy"

$ws.Range("A10").Value = 7
$ws.Range("B10").Value = 45233.65996528935
$ws.Range("C10").Value = 28
$ws.Range("D10").Value = 7
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 3
$ws.Range("G10").Value = "Fancy"

$ws.Range("A11").Value = 8
$ws.Range("B11").Value = 45233.68449534913
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = "Margerita"

$ws.Range("A12").Value = 8
$ws.Range("B12").Value = 45233.68449534913
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = 4
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = "Spicy"

$ws.Range("A13").Value = 8
$ws.Range("B13").Value = 45233.68449534913
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 7
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 3
$ws.Range("G13").Value = "Fancy"
